$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped from
# 45175 (2023-09-06) to 45177 (2023-09-08) for every data row (2..118).
$ws.Range("C2:C118").Value = 45177
